$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'323.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.27%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'39.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.85%"
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'11.77%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08031"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.98%"
$ws.Range("E5").Style = "Normal"

$ws.Range("B6").Value = "'GateToken"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "'4.587"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.28%"
$ws.Range("E6").Style = "Normal"

$ws.Range("B7").Value = "'KuCoinToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'8.671"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.14%"
$ws.Range("E7").Style = "Normal"

$ws.Range("B8").Value = "'FTXToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'1.926"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.86%"
$ws.Range("E8").Style = "Normal"

$ws.Range("B9").Value = "'BTSEToken"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'2.953"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.14%"
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = "'MXToken"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.9334"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.18%"
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.1273"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-8.39%"
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = "'WazirX"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.1969"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.07%"
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = "'MCDex"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'8.743"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'20.94%"
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = "'MandalaExchangeToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09201"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.34%"
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = "'BitrueCoin"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.03555"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'3.66%"
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = "'BitMartToken"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.1046"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'9.36%"
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = "'BitForexToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.001292"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-8.02%"
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = "'TigerCash"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'0.006151"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-3.38%"
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = "'LEO"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'3.348"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.30%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.3566"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.20%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.1419"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").Value = "'0.2448"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.69%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04412"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.69%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001261"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'3.31%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.33%"
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'-11.62%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").Value = "'0.02445"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.22%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.05250"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.28%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007424"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.63%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.009554"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.69%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.1406"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.69%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002117"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-2.45%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.009958"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'10.88%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006734"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.76%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'-0.02%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.003002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-10.16%"
$ws.Range("E48").Style = "Normal"

$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'-0.02%"
$ws.Range("E51").Style = "Normal"
